$d = $word.ActiveDocument

# Locate the paragraph that holds the "Sub:English II" subject line. It is
# the only paragraph whose text contains both "Sub:" and "English".
$paras = $d.Paragraphs
$target = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*Sub:*English*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the 'Sub:English' paragraph"
}

$range = $target.Range

# Known original WordprocessingML for this paragraph (runs "Sub:" and
# "English" sit back-to-back, flanked by stray spell/grammar proofErr
# markers left over from when the missing space between them was flagged).
$originalParaXml = '<w:p w14:paraId="5C8BAAB4" w14:textId="29F9022B" w:rsidR="00D466B6" w:rsidRDefault="00D466B6" w:rsidP="00A070E1"><w:pPr><w:pBdr><w:bottom w:val="double" w:sz="6" w:space="1" w:color="auto"/></w:pBdr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Sub:</w:t></w:r><w:r w:rsidR="004C3ED3"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>English</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="00D466B6"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="003A133E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>II</w:t></w:r><w:r w:rsidR="00123F96"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r w:rsidR="00123F96"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r w:rsidR="00123F96"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r w:rsidR="00123F96"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r w:rsidR="00123F96"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r w:rsidR="00123F96"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r w:rsidR="00123F96"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/></w:r></w:p>'

# 1) Drop the now-stale proofing-error bookmarks around "Sub:"/"English" -
#    the text they flagged is about to be fixed.
$fixedParaXml = $originalParaXml.Replace('<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>', '')
$fixedParaXml = $fixedParaXml.Replace('<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>', '')

# 2) Insert a new run containing a single space between "Sub:" and
#    "English" so the subject reads "Sub: English II" instead of
#    "Sub:English II".
$runBetween = '<w:r w:rsidR="004C3ED3"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>'
$englishRun = '<w:r w:rsidR="004C3ED3"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>English</w:t></w:r>'
$newEnglishRun = '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>English</w:t></w:r>'
$fixedParaXml = $fixedParaXml.Replace($englishRun, $runBetween + $newEnglishRun)

$xmlFragment = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$fixedParaXml
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$range.InsertXML($xmlFragment)
